$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.195153607888756
$ws.Range("C2").Value = 0.4250335635850888
$ws.Range("E2").Value = 0.2390384550666873
$ws.Range("F2").Value = 2.642434881980307
$ws.Range("G2").Value = 0.002469478123309374
$ws.Range("I2").Value = 0.9298523480855465
$ws.Range("J2").Value = 0.07674031380287261
$ws.Range("M2").Value = 0.5794862744444842
$ws.Range("N2").Value = 1.487528214287565
$ws.Range("B3").Value = 1.08599095788918
$ws.Range("C3").Value = 0.3826800574205151
$ws.Range("E3").Value = 0.2387103092931895
$ws.Range("F3").Value = 2.609556416485916
$ws.Range("G3").Value = 0.002474737504074486
$ws.Range("I3").Value = 0.9234423147897743
$ws.Range("J3").Value = 0.07722489474246785
$ws.Range("M3").Value = 0.5521185960387314
$ws.Range("N3").Value = 1.506413811459538
$ws.Range("B4").Value = 1.019516279113759
$ws.Range("C4").Value = 0.3568782838476352
$ws.Range("E4").Value = 0.2385590985326775
$ws.Range("F4").Value = 2.590954148733559
$ws.Range("G4").Value = 0.002478134767215045
$ws.Range("I4").Value = 0.920150138655174
$ws.Range("J4").Value = 0.07756929171092253
$ws.Range("M4").Value = 0.5356303733642136
$ws.Range("N4").Value = 1.518625668903734
$ws.Range("B5").Value = 0.9925645565521108
$ws.Range("C5").Value = 0.3464140118927901
$ws.Range("E5").Value = 0.2385101806680101
$ws.Range("F5").Value = 2.58377087922706
$ws.Range("G5").Value = 0.00247956156725822
$ws.Range("I5").Value = 0.9189694931319821
$ws.Range("J5").Value = 0.07772136029099386
$ws.Range("M5").Value = 0.5289904879460039
$ws.Range("N5").Value = 1.52375657140588
$ws.Range("B6").Value = 0.9880974909277711
$ws.Range("C6").Value = 0.3446794281148584
$ws.Range("E6").Value = 0.2385028268329812
$ws.Range("F6").Value = 2.582602060968753
$ws.Range("G6").Value = 0.002479801050757527
$ws.Range("I6").Value = 0.9187831422825994
$ws.Range("J6").Value = 0.07774731752411057
$ws.Range("M6").Value = 0.5278927186395563
$ws.Range("N6").Value = 1.524617870980158
$ws.Range("B7").Value = 1.01915224495923
$ws.Range("C7").Value = 0.3567369572244843
$ws.Range("E7").Value = 0.2385583872981289
$ws.Range("F7").Value = 2.590855665726892
$ws.Range("G7").Value = 0.002478153837766947
$ws.Range("I7").Value = 0.9201335656023133
$ws.Range("J7").Value = 0.07757129517030137
$ws.Range("M7").Value = 0.5355405050950779
$ws.Range("N7").Value = 1.518694241190055
$ws.Range("B8").Value = 1.157398912503879
$ws.Range("C8").Value = 0.4103872311790724
$ws.Range("E8").Value = 0.2389149121691574
$ws.Range("F8").Value = 2.630768555514194
$ws.Range("G8").Value = 0.002471256781851396
$ws.Range("I8").Value = 0.9275080402742475
$ws.Range("J8").Value = 0.07689763476517797
$ws.Range("M8").Value = 0.5699842719937322
$ws.Range("N8").Value = 1.493911852784112
$ws.Range("B9").Value = 1.432958644386758
$ws.Range("C9").Value = 0.5172603433704808
$ws.Range("E9").Value = 0.2400104589387482
$ws.Range("F9").Value = 2.72168329188267
$ws.Range("G9").Value = 0.002459057806668483
$ws.Range("I9").Value = 0.9471188295091366
$ws.Range("J9").Value = 0.07595117220711955
$ws.Range("M9").Value = 0.6400451174903594
$ws.Range("N9").Value = 1.450220288993698
$ws.Range("B10").Value = 1.638265611294912
$ws.Range("C10").Value = 0.5968737239617212
$ws.Range("E10").Value = 0.2410539280330752
$ws.Range("F10").Value = 2.796292662816455
$ws.Range("G10").Value = 0.002450894217336804
$ws.Range("I10").Value = 0.9647284474068414
$ws.Range("J10").Value = 0.07548804680949672
$ws.Range("M10").Value = 0.6930765338651668
$ws.Range("N10").Value = 1.421136146810042
$ws.Range("B11").Value = 1.732315588785696
$ws.Range("C11").Value = 0.6333470441033455
$ws.Range("E11").Value = 0.2415798017255142
$ws.Range("F11").Value = 2.831954939094089
$ws.Range("G11").Value = 0.002447351864318249
$ws.Range("I11").Value = 0.9734484945874442
$ws.Range("J11").Value = 0.07532864387620108
$ws.Range("M11").Value = 0.7175457221779169
$ws.Range("N11").Value = 1.40856535719692
$ws.Range("B12").Value = 1.768026127236624
$ws.Range("C12").Value = 0.6471967504994609
$ws.Range("E12").Value = 0.2417862420703756
$ws.Range("F12").Value = 2.845708695637626
$ws.Range("G12").Value = 0.002446034946219204
$ws.Range("I12").Value = 0.9768536349465649
$ws.Range("J12").Value = 0.07527572630507251
$ws.Range("M12").Value = 0.726861502160375
$ws.Range("N12").Value = 1.403900494790065
$ws.Range("B13").Value = 1.760330933387479
$ws.Range("C13").Value = 0.6442122579688316
$ws.Range("E13").Value = 0.2417414575510826
$ws.Range("F13").Value = 2.842735468045248
$ws.Range("G13").Value = 0.002446317481117927
$ws.Range("I13").Value = 0.9761156771883179
$ws.Range("J13").Value = 0.07528679078764
$ws.Range("M13").Value = 0.7248529623381899
$ws.Range("N13").Value = 1.404900903494649
$ws.Range("B14").Value = 1.735251585407013
$ws.Range("C14").Value = 0.6344856986647756
$ws.Range("E14").Value = 0.2415966396579243
$ws.Range("F14").Value = 2.833081464847538
$ws.Range("G14").Value = 0.002447243030563274
$ws.Range("I14").Value = 0.9737265662734842
$ws.Range("J14").Value = 0.07532414069596172
$ws.Range("M14").Value = 0.7183111370032265
$ws.Range("N14").Value = 1.408179660539453
$ws.Range("B15").Value = 1.719902294404847
$ws.Range("C15").Value = 0.6285328893791871
$ws.Range("E15").Value = 0.2415088839134221
$ws.Range("F15").Value = 2.827200613441562
$ws.Range("G15").Value = 0.002447813141888857
$ws.Range("I15").Value = 0.9722766179417732
$ws.Range("J15").Value = 0.07534799027490635
$ws.Range("M15").Value = 0.7143105763092308
$ws.Range("N15").Value = 1.410200437362292
$ws.Range("B16").Value = 1.632132544315084
$ws.Range("C16").Value = 0.5944953698332824
$ws.Range("E16").Value = 0.2410205861161678
$ws.Range("F16").Value = 2.793996837607438
$ws.Range("G16").Value = 0.00245112915334826
$ws.Range("I16").Value = 0.964172934662173
$ws.Range("J16").Value = 0.07549950204625588
$ws.Range("M16").Value = 0.6914843750933528
$ws.Range("N16").Value = 1.421971000135279
$ws.Range("B17").Value = 1.578457549946279
$ws.Range("C17").Value = 0.5736810201110529
$ws.Range("E17").Value = 0.2407341016840476
$ws.Range("F17").Value = 2.774069591008015
$ws.Range("G17").Value = 0.002453207195176988
$ws.Range("I17").Value = 0.9593840068890245
$ws.Range("J17").Value = 0.0756056329902961
$ws.Range("M17").Value = 0.6775697448406817
$ws.Range("N17").Value = 1.429361223814468
$ws.Range("B18").Value = 1.547646642996654
$ws.Range("C18").Value = 0.5617332590074398
$ws.Range("E18").Value = 0.2405741455227179
$ws.Range("F18").Value = 2.762769961275382
$ws.Range("G18").Value = 0.002454418562080288
$ws.Range("I18").Value = 0.9566962458493151
$ws.Range("J18").Value = 0.0756714986204976
$ws.Range("M18").Value = 0.6695988698243838
$ws.Range("N18").Value = 1.433673948496164
$ws.Range("B19").Value = 1.537225116222032
$ws.Range("C19").Value = 0.557692057999418
$ws.Range("E19").Value = 0.2405208171023112
$ws.Range("F19").Value = 2.758971875618641
$ws.Range("G19").Value = 0.00245483148510066
$ws.Range("I19").Value = 0.9557976408987372
$ws.Range("J19").Value = 0.07569462559584395
$ws.Range("M19").Value = 0.6669056353969864
$ws.Range("N19").Value = 1.435144805685859
$ws.Range("B20").Value = 1.584164964449712
$ws.Range("C20").Value = 0.5758942405773269
$ws.Range("E20").Value = 0.2407640997951468
$ws.Range("F20").Value = 2.776174103078091
$ws.Range("G20").Value = 0.002452984315471925
$ws.Range("I20").Value = 0.9598868853342708
$ws.Range("J20").Value = 0.0755938356052539
$ws.Range("M20").Value = 0.6790476206715397
$ws.Range("N20").Value = 1.428568094570551
$ws.Range("B21").Value = 1.742615382737824
$ws.Range("C21").Value = 0.6373415831452007
$ws.Range("E21").Value = 0.2416389784335884
$ws.Range("F21").Value = 2.83591030285757
$ws.Range("G21").Value = 0.002446970510505816
$ws.Range("I21").Value = 0.9744255009835996
$ws.Range("J21").Value = 0.07531296749559147
$ws.Range("M21").Value = 0.7202312754533438
$ws.Range("N21").Value = 1.407214015908025
$ws.Range("B22").Value = 1.846731408901007
$ws.Range("C22").Value = 0.6777233024186557
$ws.Range("E22").Value = 0.2422533052919178
$ws.Range("F22").Value = 2.876404802600405
$ws.Range("G22").Value = 0.002443182848742313
$ws.Range("I22").Value = 0.984528461969731
$ws.Range("J22").Value = 0.0751728322485441
$ws.Range("M22").Value = 0.7474377316471816
$ws.Range("N22").Value = 1.393814351739898
$ws.Range("B23").Value = 1.791110961183563
$ws.Range("C23").Value = 0.6561500994665721
$ws.Range("E23").Value = 0.2419215539073605
$ws.Range("F23").Value = 2.854658603573938
$ws.Range("G23").Value = 0.002445191383113228
$ws.Range("I23").Value = 0.9790809690140492
$ws.Range("J23").Value = 0.0752436274586934
$ws.Range("M23").Value = 0.7328904698633494
$ws.Range("N23").Value = 1.400914903048289
$ws.Range("B24").Value = 1.58158449433364
$ws.Range("C24").Value = 0.5748935854873594
$ws.Range("E24").Value = 0.2407505228570592
$ws.Range("F24").Value = 2.775222164863521
$ws.Range("G24").Value = 0.002453085027437708
$ws.Range("I24").Value = 0.9596593301831291
$ws.Range("J24").Value = 0.07559915410310936
$ws.Range("M24").Value = 0.6783793832224632
$ws.Range("N24").Value = 1.428926469039595
$ws.Range("B25").Value = 1.357921279209279
$ws.Range("C25").Value = 0.4881620077148909
$ws.Range("E25").Value = 0.2396719679884782
$ws.Range("F25").Value = 2.69572385872064
$ws.Range("G25").Value = 0.002462216960071617
$ws.Range("I25").Value = 0.9412559069706177
$ws.Range("J25").Value = 0.07616671566664124
$ws.Range("M25").Value = 0.6208199790538274
$ws.Range("N25").Value = 1.461512191497867
